$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "data as of" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Julio de 2020 a las 04:56"

# Row 44 becomes Bolivia with its updated stats
$ws.Range("A44").Value = "Bolivia"
$ws.Range("B44").Value = 33219
$ws.Range("C44").Value = 1094
$ws.Range("D44").Value = 9340
$ws.Range("E44").Value = 22756
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 52
$ws.Range("H44").Value = 1123

# Row 45 becomes Republica Dominicana, carrying what used to be row 44's data
$ws.Range("A45").Value = "Republica Dominicana"
$ws.Range("B45").Value = 32568
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 17580
$ws.Range("E45").Value = 14241
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 747

# Row 75 (Australia) data refresh
$ws.Range("B75").Value = 7920
$ws.Range("C75").Value = 84
$ws.Range("D75").Value = 7040
$ws.Range("E75").Value = 776

# Row 82 (Haiti) data refresh
$ws.Range("B82").Value = 5975
$ws.Range("C82").Value = 42
$ws.Range("E82").Value = 5031

# Rows 203-206: swap the country labels in place (Laos/Santa Lucia, Fiyi/Dominica)
$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("A204").Value = "Laos"
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"
